# Planning & Analysis.xlsx - "Start bringing in blazor samples"
#
# 1. Trim the stray note rows under the "Global Feature Analysis" table.
# 2. Insert a new "Blazor" sheet right after "Global Feature Analysis" with
#    a starter list of topics to cover.
# 3. Remove the "User Dialogs" row from "Blocked Features" (it shrinks the
#    Table6 range from A1:B5 to A1:B4).
# 4. Drop the now-unused "Framework" sheet entirely.
# 5. Restore per-sheet selections and leave "Blocked Features" as the active tab.

$wb = $excel.ActiveWorkbook

# --- 1. Global Feature Analysis: delete the loose note rows below the table ---
$ws1 = $wb.Worksheets.Item("Global Feature Analysis")
$ws1.Rows.Item(29).Delete()
$ws1.Rows.Item(27).Delete()
$ws1.Rows.Item(25).Delete()
$ws1.Rows.Item(23).Delete()
$ws1.Rows.Item(22).Delete()
$ws1.Rows.Item(21).Delete()
$ws1.Rows.Item(19).Delete()
$ws1.Range("D18").Select()

# --- 2. Insert the new "Blazor" sheet after "Global Feature Analysis" ---
$blazor = $wb.Worksheets.Add($null, $ws1)
$blazor.Name = "Blazor"

# Populate in the same first-seen order the authoring session used (matches
# the resulting shared-string table), even though the rows land as A1..A10.
$blazor.Cells.Item(2, 1).Value = "Push"
$blazor.Cells.Item(3, 1).Value = "Notifications"
$blazor.Cells.Item(1, 1).Value = "BluetoothLE"
$blazor.Cells.Item(4, 1).Value = "Locations - GPS only"
$blazor.Cells.Item(5, 1).Value = "Sensors"
$blazor.Cells.Item(6, 1).Value = "Core - Connectivity"
$blazor.Cells.Item(7, 1).Value = "Core - Battery"
$blazor.Cells.Item(8, 1).Value = "Jobs"
$blazor.Cells.Item(9, 1).Value = "Speech Recognition"
$blazor.Cells.Item(10, 1).Value = "NFC"
$blazor.Range("A3").Select()

# --- 3. Blocked Features: drop the "User Dialogs" row ---
$blocked = $wb.Worksheets.Item("Blocked Features")
$blocked.Rows.Item(2).Delete()

# --- 4. Remove the "Framework" sheet ---
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Framework").Delete()
$excel.DisplayAlerts = $true

# --- 5. Restore selections on the untouched sheets, activate Blocked Features ---
$platform = $wb.Worksheets.Item("Platform")
$platform.Range("B12").Select()

$featureMovement = $wb.Worksheets.Item("Feature Movement")
$featureMovement.Range("B13").Select()

$blocked.Select()
$blocked.Range("B3").Select()
